# Apply cell-value updates for the cryptos list (price + 1h volume change refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.473.60'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.625.21'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.58'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.18'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0854'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '1.855.24'
$ws.Range('D13').Value = '1.629.91'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.96'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '234.84'
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('D18').Value = '26.500.47'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '0.0₃0726'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.30'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.77'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').Value = '1.516.97'
$ws.Range('E32').Value = '  +4.94%  '
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.568'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0166'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.834'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.86'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').Value = '1.765.66'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.85'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.762'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('E46').Value = '  -2.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.83'
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0501'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.57'
$ws.Range('E51').Value = '  +1.12%  '
